$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.421.17'
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").Value = '1.830.33'
$ws.Range("E3").Value = '  -1.52%  '
$ws.Range("E4").Value = '  -3.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.51%  '
$ws.Range("E6").Value = '  -2.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4302'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3705'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07254'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8677'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '1.837.93'
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.687'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.363'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07069'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.92'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008913'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.18%  '
$ws.Range("E19").Value = '  -2.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.01%  '
$ws.Range("D21").Value = '27.430.68'
$ws.Range("E21").Value = '  -1.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.174'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.58%  '
$ws.Range("D24").Value = '2.062.57'
$ws.Range("E24").Value = '  -1.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.018'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.45'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.140'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.302'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.55'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08864'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.210'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.66%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7697'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.507'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.899'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.17%  '
$ws.Range("E36").Value = '  -3.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.123'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01965'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05286'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.170'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.880'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1679'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5083'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.692'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.66'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '106.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4740'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06422'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.60%  '
$ws.Range("E49").Value = '  -3.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.674'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.830'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.18%  '
